# Applies the "Revert "Revert ...""" changes described by the diff:
#  - Slide 1, shape "직사각형 9": "Win 32 API" -> "Utilities", plus a new
#    paragraph "(Include Win 32 API)" appended underneath it.
#  - Slide 3, shapes "직사각형 56" / "직사각형 57": merge the two-run
#    "Page "+"Setter" / "Screen "+"Getter" texts into single runs.
#  - Slide 5, shape "직사각형 26": "Win 32 API" -> "Utilities" (simple swap).
#  - Slide 5, shapes "직사각형 30" / "직사각형 31": same run-merge as slide 3.
#
# Note: when the new text, once the runs are concatenated, is identical to
# the text already displayed (e.g. merging "Page " + "Setter" into a
# single run "Page Setter" - same visible characters, just fewer runs) the
# host's TextRange.Text setter treats it as a no-op and leaves the original
# multi-run split alone. Routing the assignment through a throwaway value
# first forces the host to actually rebuild the run(s), after which we set
# the real text.

$p = $ppt.ActivePresentation

# --- Slide 1: "Win 32 API" -> "Utilities" + new "(Include Win 32 API)" para ---
$slide1 = $p.Slides.Item(1)
$shp1 = $slide1.Shapes.Item(6)
$tr1 = $shp1.TextFrame.TextRange
$tr1.Text = "Utilities"
$tr1.InsertAfter("`r(Include Win 32 API)") | Out-Null

# --- Slide 3: merge "Page " + "Setter" runs, and "Screen " + "Getter" runs ---
$slide3 = $p.Slides.Item(3)

$tr3a = $slide3.Shapes.Item(11).TextFrame.TextRange
$tr3a.Text = "_tmp_"
$tr3a.Text = "Page Setter"

$tr3b = $slide3.Shapes.Item(12).TextFrame.TextRange
$tr3b.Text = "_tmp_"
$tr3b.Text = "Screen Getter"

# --- Slide 5: "Win 32 API" -> "Utilities" (no extra paragraph this time) ---
$slide5 = $p.Slides.Item(5)
$slide5.Shapes.Item(7).TextFrame.TextRange.Text = "Utilities"

# --- Slide 5: merge "Page " + "Setter" runs, and "Screen " + "Getter" runs ---
$tr5a = $slide5.Shapes.Item(11).TextFrame.TextRange
$tr5a.Text = "_tmp_"
$tr5a.Text = "Page Setter"

$tr5b = $slide5.Shapes.Item(12).TextFrame.TextRange
$tr5b.Text = "_tmp_"
$tr5b.Text = "Screen Getter"
